# Auto-generated updates for cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PEPE row (29) price contains a Unicode subscript-3 character (U+2083)
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.604.77"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.624.82"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.02"
$ws.Range("E5").Value = "  +3.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.62"
$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.54"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.376"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.093.02"
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.31"
$ws.Range("E14").Value = "  +12.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.604.17"
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.638.03"
$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.55"
$ws.Range("E18").Value = "  +2.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.71"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.46"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.92"
$ws.Range("E21").Value = "  -1.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.526"
$ws.Range("E23").Value = "  -1.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.63"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.160"
$ws.Range("E26").Value = "  +1.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.19"
$ws.Range("E27").Value = "  +6.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.05"
$ws.Range("E28").Value = "  +12.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = ("0.0{0}0804" -f $sub3)
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.50"
$ws.Range("E30").Value = "  +2.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.37"
$ws.Range("E31").Value = "  +5.49%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("E34").Value = "  +4.31%  "

$ws.Range("E35").Value = "  +5.31%  "

$ws.Range("E36").Value = "  +7.92%  "

$ws.Range("E37").Value = "  +3.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "332.51"
$ws.Range("E38").Value = "  +12.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.02"
$ws.Range("E39").Value = "  +5.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.78"
$ws.Range("E40").Value = "  +2.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.857"
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.18"
$ws.Range("E42").Value = "  +6.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "134.71"
$ws.Range("E43").Value = "  -4.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.11"
$ws.Range("E44").Value = "  +2.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0995"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("E47").Value = "  +0.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0555"
$ws.Range("E48").Value = "  +1.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.27"
$ws.Range("E49").Value = "  +3.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0244"
$ws.Range("E50").Value = "  +1.90%  "

$ws.Range("E51").Value = "  +0.59%  "
